$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'Saco Para Até 10 Bolas / Porta Bolas Fio 4mm Branco'
$ws.Range('B2').Value = 12.81
$ws.Range('C2').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_720045-MLU72700298625_112023-E.webp'
$ws.Range('D2').Value = 'https://www.mercadolivre.com.br/saco-para-ate-10-bolas-porta-bolas-fio-4mm-branco/p/MLB25996601#polycard_client=search-nordic&searchVariation=MLB25996601&wid=MLB4147585512&position=5&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A3').Value = 'Bola De Futebol Para Campo Tamanho 5 Petrin Presente Cor'
$ws.Range('B3').Value = 18.37
$ws.Range('C3').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_859381-MLA81212900247_122024-E.webp'
$ws.Range('D3').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-para-campo-tamanho-5-petrin-presente-cor/p/MLB44556814#polycard_client=search-nordic&searchVariation=MLB44556814&wid=MLB4007228671&position=39&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A4').Value = 'Trave Futebol Gol Kit Mini Golzinho Infantil Brinquedo Bola'
$ws.Range('B4').Value = 33.54
$ws.Range('C4').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_764435-MLA79889067517_102024-E.webp'
$ws.Range('D4').Value = 'https://www.mercadolivre.com.br/trave-futebol-gol-kit-mini-golzinho-infantil-brinquedo-bola/p/MLB41775954#polycard_client=search-nordic&searchVariation=MLB41775954&wid=MLB5118193404&position=11&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A5').Value = 'Tabela De Medição De Pressão, Medidor De Ponteiro De Futebol'
$ws.Range('B5').Value = 38.31
$ws.Range('C5').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_725638-CBT81594394892_012025-E.webp'
$ws.Range('D5').Value = 'https://www.mercadolivre.com.br/pressure-measurement-table-portable-football-pointer-gauge/p/MLB2009472856#polycard_client=search-nordic&searchVariation=MLB2009472856&wid=MLB4000139541&position=9&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A6').Value = 'Bola Futebol Tamanho 5 Capotão + Bomba De Ar E Bico Cor Azul'
$ws.Range('B6').Value = 39.99
$ws.Range('C6').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_628564-MLU78026959026_082024-E.webp'
$ws.Range('D6').Value = 'https://www.mercadolivre.com.br/bola-futebol-tamanho-5-capoto-bomba-de-ar-e-bico-cor-azul/p/MLB39157733#polycard_client=search-nordic&searchVariation=MLB39157733&wid=MLB3870171799&position=30&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A7').Value = 'Kit Bola De Vôlei + Bola De Futebol + Bomba Pró Sortidas'
$ws.Range('B7').Value = 41.22
$ws.Range('C7').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_702577-MLA82266657438_022025-E.webp'
$ws.Range('D7').Value = 'https://www.mercadolivre.com.br/kit-bola-de-vlei-bola-de-futebol-bomba-pro-sortidas/p/MLB46184165#polycard_client=search-nordic&searchVariation=MLB46184165&wid=MLB5292902016&position=22&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A8').Value = 'Bola De Vôlei 3.5 Texturizada Amarela Rainha Topper'
$ws.Range('B8').Value = 44.99
$ws.Range('C8').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_875461-MLU78119836106_082024-E.webp'
$ws.Range('D8').Value = 'https://www.mercadolivre.com.br/bola-de-vlei-35-texturizada-amarela-rainha-topper/p/MLB21012198#polycard_client=search-nordic&searchVariation=MLB21012198&wid=MLB3963193225&position=27&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A9').Value = 'Bola De Futebol Luke Sports Tamanho 5 Campo E Society Branca Cor Branco'
$ws.Range('B9').Value = 48.29
$ws.Range('C9').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_921199-MLA80927520322_122024-E.webp'
$ws.Range('D9').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-luke-sports-tamanho-5-campo-e-society-branca-cor-branco/p/MLB44516604#polycard_client=search-nordic&searchVariation=MLB44516604&wid=MLB5214268684&position=32&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A10').Value = '10 Bolas De Vinil Futebol Dente De Leite C/ Cor Branco Com Preto'
$ws.Range('B10').Value = 48.97
$ws.Range('C10').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_785671-MLU77640619462_072024-E.webp'
$ws.Range('D10').Value = 'https://www.mercadolivre.com.br/10-bolas-de-vinil-futebol-dente-de-leite-c-cor-branco-com-preto/p/MLB38350019#polycard_client=search-nordic&searchVariation=MLB38350019&wid=MLB3855446145&position=24&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A11').Value = 'Bola De Futebol De Campo Slick 2020 Topper Cor Amarelo Neon/Preto'
$ws.Range('B11').Value = 48.99
$ws.Range('C11').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_969183-MLU75988202853_042024-E.webp'
$ws.Range('D11').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-de-campo-slick-2020-topper-cor-amarelo-neonpreto/p/MLB19754376#polycard_client=search-nordic&searchVariation=MLB19754376&wid=MLB5277328830&position=38&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A12').Value = 'Bola De Futebol Slick Campo Cor Azul/Branco/Preto Topper'
$ws.Range('B12').Value = 50.72
$ws.Range('C12').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_979047-MLU71168535619_082023-E.webp'
$ws.Range('D12').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-slick-campo-cor-azulbrancopreto-topper/p/MLB19765484#polycard_client=search-nordic&searchVariation=MLB19765484&wid=MLB5270990860&position=35&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A13').Value = 'Bola Futebol Society Ótimo Custo Beneficio Peso Ideal Top Cor Azul'
$ws.Range('B13').Value = 53.82
$ws.Range('C13').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_862896-MLU69453002841_052023-E.webp'
$ws.Range('D13').Value = 'https://www.mercadolivre.com.br/bola-futebol-society-otimo-custo-beneficio-peso-ideal-top-cor-azul/p/MLB23339357#polycard_client=search-nordic&searchVariation=MLB23339357&wid=MLB3923038695&position=6&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A14').Value = 'Bola De Futebol de Campo Slick 24 Tecnofusion Multicor Topper'
$ws.Range('B14').Value = 53.99
$ws.Range('C14').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_722050-MLU76926100481_062024-E.webp'
$ws.Range('D14').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-de-campo-slick-24-tecnofusion-multicor-topper/p/MLB35378541#polycard_client=search-nordic&searchVariation=MLB35378541&wid=MLB3963108755&position=40&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A15').Value = 'Bola De Futebol De Campo Slick 2020 Topper'
$ws.Range('B15').Value = 54.99
$ws.Range('C15').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_735756-MLB54078855256_022023-E.webp'
$ws.Range('D15').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-de-campo-slick-2020-topper/p/MLB22262774#polycard_client=search-nordic&searchVariation=MLB22262774&wid=MLB5270990842&position=13&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A16').Value = 'Bola De Futsal Slick 24 Câmara Airvility 6 Gomos Topper'
$ws.Range('B16').Value = 54.99
$ws.Range('C16').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_790020-MLU75616505649_042024-E.webp'
$ws.Range('D16').Value = 'https://www.mercadolivre.com.br/bola-de-futsal-slick-24-cmara-airvility-6-gomos-topper/p/MLB35378564#polycard_client=search-nordic&searchVariation=MLB35378564&wid=MLB3963063089&position=29&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A17').Value = 'Bola De Futebol Para Campo 22 I Topper Cor Branco/Verde Mar/Preto'
$ws.Range('B17').Value = 55.85
$ws.Range('C17').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_967540-MLU69495823375_052023-E.webp'
$ws.Range('D17').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-para-campo-22-i-topper-cor-brancoverde-marpreto/p/MLB22272252#polycard_client=search-nordic&searchVariation=MLB22272252&wid=MLB4198884500&position=4&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A18').Value = 'Bola De Futebol Campo Slick Cup Topper'
$ws.Range('B18').Value = 59.15
$ws.Range('C18').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_659913-MLA80298884457_102024-E.webp'
$ws.Range('D18').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-campo-slick-cup-topper/p/MLB19698380#polycard_client=search-nordic&searchVariation=MLB19698380&wid=MLB5189140588&position=8&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A19').Value = 'Bola De Futsal Slick 22 Topper Cor Branco/Amarelo Neon/Azul'
$ws.Range('B19').Value = 59.2
$ws.Range('C19').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_656329-MLB54078855245_022023-E.webp'
$ws.Range('D19').Value = 'https://www.mercadolivre.com.br/bola-de-futsal-slick-22-topper-cor-brancoamarelo-neonazul/p/MLB22262771#polycard_client=search-nordic&searchVariation=MLB22262771&wid=MLB3820042822&position=12&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A20').Value = 'Bola De Futsal Slick Cup Topper Cor Laranja/Azul Anodizado/Preto'
$ws.Range('B20').Value = 59.9
$ws.Range('C20').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_882043-MLU75981819389_042024-E.webp'
$ws.Range('D20').Value = 'https://www.mercadolivre.com.br/bola-de-futsal-slick-cup-topper-cor-laranjaazul-anodizadopreto/p/MLB19710248#polycard_client=search-nordic&searchVariation=MLB19710248&wid=MLB2889188605&position=36&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A21').Value = 'Topper Futsal Slick Bola Futebol Salão Branco Com Azul'
$ws.Range('B21').Value = 60.99
$ws.Range('C21').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_880504-MLU71134471481_082023-E.webp'
$ws.Range('D21').Value = 'https://www.mercadolivre.com.br/topper-futsal-slick-bola-futebol-salo-branco-com-azul/p/MLB26148855#polycard_client=search-nordic&searchVariation=MLB26148855&wid=MLB3963093661&position=17&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A22').Value = 'Bola De Futebol Society Slick 2020 Cor Branco/Preto/Amarelo Tamanho Único Topper'
$ws.Range('B22').Value = 61.76
$ws.Range('C22').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_869053-MLB51937100243_102022-E.webp'
$ws.Range('D22').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-society-slick-2020-cor-brancopretoamarelo-tamanho-unico-topper/p/MLB19744266#polycard_client=search-nordic&searchVariation=MLB19744266&wid=MLB5222132312&position=20&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A23').Value = 'Bola De Futevolei Altinha N°5 Quadra Areia Pro Praia'
$ws.Range('B23').Value = 65.97
$ws.Range('C23').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_645076-MLB82458925554_022025-E-bola-de-futevolei-altinha-n5-quadra-areia-pro-praia.webp'
$ws.Range('D23').Value = 'https://click1.mercadolivre.com.br/mclics/clicks/external/MLB/count?a=mYpbkfZv%2BWrkHXFcJzH%2BLNwYdTR0CClKwqUNkB8gCyWC%2FkdDlw6z8okZwVrmsuyR2mvs0GYH0ptERW22b6tFsUz0f5KJyKt0nZlkml5NHIihom1Px0J1CrbZ6BSH%2FQTGQE6B9ax%2B8GkDXee%2B3gEdsQRcITpBPk4s%2Fw8gpWzkpb6nHYDnDCDl2xtI3hK1lB3Wd5%2F%2Fd%2F80G3ATney9LcSZw2H4xE2MA85NvXkt7eZqjfydIIsRHfEjNqfi%2FkCWbCW44xypAzY70XuE9F34Sb0BBB5IrTcC1McHT%2FQt4QIlLn5MravpRWVd4SjFa6QU9%2FCOQJCek3J%2BfhFmi6pNVPwdjvRnCYAj5w617jbLWvuC9BJT0I8D7bMaV8QDGZALYfTZOfGT97MfuUXaR8ghh0ObguNGxU55WI2FSBFhZdflUaZsCNlQoxpEcAePwzEG3z85QfUHPHfZqR5w0na0JYqMbdr4SD8yNvtmACVwSDGRi44oCy3mDQxN6SD%2B0ksm4DoC%2BV3xPf2OEQcAkZByaZrXT4%2B3gafEjiieOk%2FK3hq79iySWUkvbfQ9s%2FKjQRQuB2SFp7BdKa3JXs9ckJwlOhqDwF9W6qK6NtExi5tDMVlig4EIBE85sqQAg46ddNkNlAHoWDGbafgKjwSTAxNppif8QmfMEE%2B26ocAjnPEHS8rDpgCfzEUyAQrou8F0S38ZvEluacn6hROTtflt0mPDQE7axmgsZnO2W9VsgeQreDZnbUF1SNZwEQeM11LfwQqyyGH%2Bme8eNvc8LGIA8UlFpcKtslbO%2FoJPAGu4A%2B0tShxXMRtpBbPDC9usKwaFGuPT7W5MhE6hmd3jfwtV84NEsMGcnY5tzYenexutjdg&rb=x#polycard_client=search-nordic&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Range('A24').Value = 'Bola De Futevolei Altinha N°5 Quadra Areia Pro Praia'
$ws.Range('B24').Value = 65.97
$ws.Range('C24').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_645076-MLB82458925554_022025-E-bola-de-futevolei-altinha-n5-quadra-areia-pro-praia.webp'
$ws.Range('D24').Value = 'https://produto.mercadolivre.com.br/MLB-4006687617-bola-de-futevolei-altinha-n5-quadra-areia-pro-praia-_JM?searchVariation=183156389516#polycard_client=search-nordic&searchVariation=183156389516&position=44&search_layout=grid&type=item&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Range('A25').Value = 'Bola De Futevôlei Diadora Elite-r'
$ws.Range('B25').Value = 69.59999999999999
$ws.Range('C25').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_777673-MLB74305800482_022024-E-bola-de-futevlei-diadora-elite-r.webp'
$ws.Range('D25').Value = 'https://click1.mercadolivre.com.br/mclics/clicks/external/MLB/count?a=FA2v80kO0Q62jIC3APOe8OLdHZosDSTo1hLLRyRzr90a4JOmG9fNJuYmnoO8wg22iqTmr27d%2BrK3sZOAMDcUxTXULya%2FuqG3Uce7FwirY265ddOp%2FaRo0%2FUfe1%2Fe8RReqaCkhX4gJYzGKk62Jq9YuTkU%2FZbzoSwUeYEBMevSBPog0Bu0gIxXMJ93cIk4gUOvWjnzprMXLSMEiSakhNXOETkZc5KUFk1bQXqKVFLpHWMuz2XLh68FPbRohnW7tl%2BCcmfgQvwQKhoklCxqKO8wcN%2B12x61YphpgBIMTsCsHnlUsyjm18ZhifBDehyGu9uKFLf5gpt612eYaSwpb3psZ42tp7JY20HPiQBPpcWki4%2FqHPJIE8iZwN9PD8PRPxmFi%2B%2F9jhfUYSfXKW7I6emq3gWPDPAro%2BtAPq1Mlbu1WvjcEOFDOeYK%2F4o56stlpSmFH4zTluHTF8HsKZqH49c8i6VjRDwLf4JiWGWPeyy5JuASJg7cJfmPAASSftaG%2B5b%2BFAQx5qMBbAgkFrIVkMakbpJyEXx8LKisH3WrSN3mVVhT%2Bal25Bev9FVuqCIkH5rpVvnIOIthHtGurHdG5%2FUMhpGuYitAJxoXDha%2FceYFueTm8Ya0XdxcZd4LtZXVUsnvCxj%2FjN%2B5R%2BqrYb4uyMOPasbyNJqMjg4nxO6hb6fRGkW5imtNH17dsnmSd22hBtwdRfUtVo6b7VeYEEiD72oR%2Bf2sqpC8IVbMa%2B5uV03kiM3O7Atg3Dc6%2FLHKx%2Blem1khGGw4rvQ51VV%2FLhtFT0uDGeA367TVDZlCRBGdGq%2F8EF4H0r0oG9GbdOGxqt7Ld0bygm%2FBhLfnOQ%3D%3D&rb=x#polycard_client=search-nordic&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Range('A26').Value = 'Bola De Futevôlei Diadora Elite-r'
$ws.Range('B26').Value = 69.59999999999999
$ws.Range('C26').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_777673-MLB74305800482_022024-E-bola-de-futevlei-diadora-elite-r.webp'
$ws.Range('D26').Value = 'https://produto.mercadolivre.com.br/MLB-3705513048-bola-de-futevlei-diadora-elite-r-_JM?searchVariation=177953619994#polycard_client=search-nordic&searchVariation=177953619994&position=50&search_layout=grid&type=item&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Range('A27').Value = 'Bola De Futevôlei Diadora Elite-r Cor Amarelo'
$ws.Range('B27').Value = 74.92
$ws.Range('C27').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_699842-MLU70081962810_062023-E.webp'
$ws.Range('D27').Value = 'https://www.mercadolivre.com.br/bola-de-futevlei-diadora-elite-r-cor-amarelo/p/MLB21096382#polycard_client=search-nordic&searchVariation=MLB21096382&wid=MLB3383487155&position=19&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A28').Value = 'Bola De Campo Topper Slick Gold Dourada Oficial Bomba De Ar'
$ws.Range('B28').Value = 74.95
$ws.Range('C28').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_869593-MLB82719894389_022025-E-bola-de-campo-topper-slick-gold-dourada-oficial-bomba-de-ar.webp'
$ws.Range('D28').Value = 'https://produto.mercadolivre.com.br/MLB-3987100373-bola-de-campo-topper-slick-gold-dourada-oficial-bomba-de-ar-_JM?searchVariation=183000579442#polycard_client=search-nordic&searchVariation=183000579442&position=43&search_layout=grid&type=item&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Range('A29').Value = 'Bola Para Futebol De Campo Bravo Xxiv Branco e Azul Penalty'
$ws.Range('B29').Value = 77.90000000000001
$ws.Range('C29').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_698520-MLU75981724023_042024-E.webp'
$ws.Range('D29').Value = 'https://www.mercadolivre.com.br/bola-para-futebol-de-campo-bravo-xxiv-branco-e-azul-penalty/p/MLB27895686#polycard_client=search-nordic&searchVariation=MLB27895686&wid=MLB3951567371&position=3&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A30').Value = 'Bola Futsal Penalty Rx 500 Xxiii Cor Azul'
$ws.Range('B30').Value = 78.90000000000001
$ws.Range('C30').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_722195-MLU77108763060_062024-E.webp'
$ws.Range('D30').Value = 'https://www.mercadolivre.com.br/bola-futsal-penalty-rx-500-xxiii-cor-azul/p/MLB25820360#polycard_client=search-nordic&searchVariation=MLB25820360&wid=MLB3996965885&position=16&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A31').Value = 'Bola De Futebol Branco E Laranja De Campo Bravo Xxiv Penalty'
$ws.Range('B31').Value = 78.98999999999999
$ws.Range('C31').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_764251-MLU75813169222_042024-E.webp'
$ws.Range('D31').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-branco-e-laranja-de-campo-bravo-xxiv-penalty/p/MLB27893082#polycard_client=search-nordic&searchVariation=MLB27893082&wid=MLB5056295770&position=23&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A32').Value = 'Bola Para Futebol De Campo Líder Xxiv Cor Branco/Azul/Verde Penalty'
$ws.Range('B32').Value = 99.98999999999999
$ws.Range('C32').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_895518-MLU74287171999_012024-E.webp'
$ws.Range('D32').Value = 'https://www.mercadolivre.com.br/bola-para-futebol-de-campo-lider-xxiv-cor-brancoazulverde-penalty/p/MLB30924589#polycard_client=search-nordic&searchVariation=MLB30924589&wid=MLB3937962097&position=10&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A33').Value = 'Bola Allpha Futsal Pró+ Madri Pu Maciez Treino Original'
$ws.Range('B33').Value = 106.69
$ws.Range('C33').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_965545-MLB82741758954_032025-E-bola-allpha-futsal-pro-madri-pu-maciez-treino-original.webp'
$ws.Range('D33').Value = 'https://produto.mercadolivre.com.br/MLB-4000309717-bola-allpha-futsal-pro-madri-pu-maciez-treino-original-_JM?searchVariation=183112590022#polycard_client=search-nordic&searchVariation=183112590022&position=45&search_layout=grid&type=item&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Range('A34').Value = 'Bola de Futebol Penalty Society Lider XXIV Branco Laranja Azul 5'
$ws.Range('B34').Value = 111.71
$ws.Range('C34').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_903645-MLU74833063908_032024-E.webp'
$ws.Range('D34').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-penalty-society-lider-xxiv-branco-laranja-azul-5/p/MLB27892495#polycard_client=search-nordic&searchVariation=MLB27892495&wid=MLB3951591553&position=18&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A35').Value = 'Bola Society Líder Xxiv Pu Laminado 6 Gomos Penalty'
$ws.Range('B35').Value = 114.99
$ws.Range('C35').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_609551-MLU78177695411_082024-E.webp'
$ws.Range('D35').Value = 'https://www.mercadolivre.com.br/bola-society-lider-xxiv-pu-laminado-6-gomos-penalty/p/MLB27924191#polycard_client=search-nordic&searchVariation=MLB27924191&wid=MLB3938335813&position=7&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A36').Value = 'Bola De Futebol Corinthians Time Oficial Licenciada Timão'
$ws.Range('B36').Value = 119
$ws.Range('C36').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_707361-MLB82839791767_032025-E-bola-de-futebol-corinthians-time-oficial-licenciada-timo.webp'
$ws.Range('D36').Value = 'https://produto.mercadolivre.com.br/MLB-3992563623-bola-de-futebol-corinthians-time-oficial-licenciada-timo-_JM?searchVariation=187103662249#polycard_client=search-nordic&searchVariation=187103662249&position=49&search_layout=grid&type=item&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Range('A37').Value = 'Bola Futevolei Altinha Ftv Ft5 Penta Futemesa Rio Cor Preto'
$ws.Range('B37').Value = 126.85
$ws.Range('C37').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_764380-MLU73463308040_122023-E.webp'
$ws.Range('D37').Value = 'https://www.mercadolivre.com.br/bola-futevolei-altinha-ftv-ft5-penta-futemesa-rio-cor-preto/p/MLB29123039#polycard_client=search-nordic&searchVariation=MLB29123039&wid=MLB4321183772&position=28&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A38').Value = 'Bola 81 Dalponte Star Campo Costurada A Mão Original'
$ws.Range('B38').Value = 131
$ws.Range('C38').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_947113-MLB45202965881_032021-E.webp'
$ws.Range('D38').Value = 'https://www.mercadolivre.com.br/bola-81-dalponte-star-campo-costurada-a-mo-original/p/MLB20563342#polycard_client=search-nordic&searchVariation=MLB20563342&wid=MLB4006264629&position=26&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A39').Value = 'Bola de futevolei Poker Rio Poker nº 5 Unidade x 1 unidades cor rosa e preto'
$ws.Range('B39').Value = 135
$ws.Range('C39').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_904045-MLU79036618824_092024-E.webp'
$ws.Range('D39').Value = 'https://www.mercadolivre.com.br/bola-de-futevolei-poker-rio-poker-n-5-unidade-x-1-unidades-cor-rosa-e-preto/p/MLB19698706#polycard_client=search-nordic&searchVariation=MLB19698706&wid=MLB5051108460&position=34&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A40').Value = 'Bola de Futevolei Poker Rio Poker nº 5 cor amarelo e preto'
$ws.Range('B40').Value = 139.29
$ws.Range('C40').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_753390-MLU74696438844_022024-E.webp'
$ws.Range('D40').Value = 'https://www.mercadolivre.com.br/bola-de-futevolei-poker-rio-poker-n-5-cor-amarelo-e-preto/p/MLB19698709#polycard_client=search-nordic&searchVariation=MLB19698709&wid=MLB3889075523&position=31&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A41').Value = 'Bola de futevolei Poker Rio Poker nº 5 Unidade x 1 unidades cor branco e preto'
$ws.Range('B41').Value = 140.11
$ws.Range('C41').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_773647-MLA79588973813_092024-E.webp'
$ws.Range('D41').Value = 'https://www.mercadolivre.com.br/bola-de-futevolei-poker-rio-poker-n-5-unidade-x-1-unidades-cor-branco-e-preto/p/MLB19698708#polycard_client=search-nordic&searchVariation=MLB19698708&wid=MLB3889079001&position=33&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A42').Value = 'Bola Futevôlei Vulcanizada Pu Poker Original Branca/azul'
$ws.Range('B42').Value = 149.9
$ws.Range('C42').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_601980-MLB83086122841_032025-E-bola-futevlei-vulcanizada-pu-poker-original-brancaazul.webp'
$ws.Range('D42').Value = 'https://produto.mercadolivre.com.br/MLB-5318075606-bola-futevlei-vulcanizada-pu-poker-original-brancaazul-_JM#polycard_client=search-nordic&position=48&search_layout=grid&type=item&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Range('A43').Value = 'Bola De Futebol Society 8 X Penalty'
$ws.Range('B43').Value = 154.71
$ws.Range('C43').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_818310-MLU76366734942_052024-E.webp'
$ws.Range('D43').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-society-8-x-penalty/p/MLB19776132#polycard_client=search-nordic&searchVariation=MLB19776132&wid=MLB3829143887&position=15&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A44').Value = 'Bola De Futebol Campo Uhlsport Match R1 Brasileirão 2025'
$ws.Range('B44').Value = 155.99
$ws.Range('C44').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_857447-MLB82813326112_032025-E-bola-de-futebol-campo-uhlsport-match-r1-brasileiro-2025.webp'
$ws.Range('D44').Value = 'https://produto.mercadolivre.com.br/MLB-4003764129-bola-de-futebol-campo-uhlsport-match-r1-brasileiro-2025-_JM#polycard_client=search-nordic&position=46&search_layout=grid&type=item&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Range('A45').Value = 'Bola De Campo Penalty S11 R2 Xxv ( Camp.paulista)'
$ws.Range('B45').Value = 159.9
$ws.Range('C45').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_687100-MLB82940622023_032025-E-bola-de-campo-penalty-s11-r2-xxv-camppaulista.webp'
$ws.Range('D45').Value = 'https://produto.mercadolivre.com.br/MLB-3996246671-bola-de-campo-penalty-s11-r2-xxv-camppaulista-_JM?searchVariation=187147572233#polycard_client=search-nordic&searchVariation=187147572233&position=47&search_layout=grid&type=item&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Range('A46').Value = 'Bola Futebol De Campo Bola 8 Penalty Cor Preto'
$ws.Range('B46').Value = 164
$ws.Range('C46').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_888107-MLA81942694869_012025-E.webp'
$ws.Range('D46').Value = 'https://www.mercadolivre.com.br/bola-futebol-de-campo-bola-8-penalty-cor-preto/p/MLB22210458#polycard_client=search-nordic&searchVariation=MLB22210458&wid=MLB3945165169&position=37&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A47').Value = 'Bola De Futebol Dalponte 81 Nitro Society Costurada À Mão Cor Branco'
$ws.Range('B47').Value = 176.62
$ws.Range('C47').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_756582-MLU74393181908_022024-E.webp'
$ws.Range('D47').Value = 'https://www.mercadolivre.com.br/bola-de-futebol-dalponte-81-nitro-society-costurada-mo-cor-branco/p/MLB32946547#polycard_client=search-nordic&searchVariation=MLB32946547&wid=MLB3596272369&position=21&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A48').Value = 'Bola Futebol Society 8 Pro Xxiv Tamanho 5 Cor Preto e Branco Penalty'
$ws.Range('B48').Value = 297
$ws.Range('C48').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_819669-MLU75482475121_032024-E.webp'
$ws.Range('D48').Value = 'https://www.mercadolivre.com.br/bola-futebol-society-8-pro-xxiv-tamanho-5-cor-preto-e-branco-penalty/p/MLB35142675#polycard_client=search-nordic&searchVariation=MLB35142675&wid=MLB4630087278&position=14&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A49').Value = 'Bola De Futsal Penalty Max 1000 Xxiv Cor Azul'
$ws.Range('B49').Value = 298.28
$ws.Range('C49').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_973308-MLU75866936168_042024-E.webp'
$ws.Range('D49').Value = 'https://www.mercadolivre.com.br/bola-de-futsal-penalty-max-1000-xxiv-cor-azul/p/MLB36227619#polycard_client=search-nordic&searchVariation=MLB36227619&wid=MLB3937932227&position=25&search_layout=grid&type=product&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d&sid=search'

$ws.Range('A50').Value = 'Bola Uhlsport Campo Game Pro Brasileirão Série C E D 2025'
$ws.Range('B50').Value = 379
$ws.Range('C50').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_652481-MLB82865558482_032025-E-bola-uhlsport-campo-game-pro-brasileiro-serie-c-e-d-2025.webp'
$ws.Range('D50').Value = 'https://produto.mercadolivre.com.br/MLB-4006258699-bola-uhlsport-campo-game-pro-brasileiro-serie-c-e-d-2025-_JM?searchVariation=187252357205#polycard_client=search-nordic&searchVariation=187252357205&position=42&search_layout=grid&type=item&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Range('A51').Value = 'Kit 50 Bolas Futebol Couro Sintético Costurada N°05'
$ws.Range('B51').Value = 916
$ws.Range('C51').Value = 'https://http2.mlstatic.com/D_Q_NP_2X_871334-MLB82990085837_032025-E-kit-50-bolas-futebol-couro-sintetico-costurada-n05.webp'
$ws.Range('D51').Value = 'https://produto.mercadolivre.com.br/MLB-5313103164-kit-50-bolas-futebol-couro-sintetico-costurada-n05-_JM?searchVariation=183096595076#polycard_client=search-nordic&searchVariation=183096595076&position=41&search_layout=grid&type=item&tracking_id=8e611849-75a7-4304-8cbe-30b2b6ffa88d'

$ws.Rows.Item(53).Delete()
$ws.Rows.Item(52).Delete()
